$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.367.38'
$ws.Range("E2").Value = '  +1.69%  '

$ws.Range("D3").Value = '1.832.79'
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("E4").Value = '  +0.88%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.53'
$ws.Range("E5").Value = '  +1.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("E6").Value = '  +0.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4741'
$ws.Range("E7").Value = '  +2.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3694'
$ws.Range("E8").Value = '  +1.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07465'
$ws.Range("E9").Value = '  +1.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8863'
$ws.Range("E10").Value = '  +2.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.45'
$ws.Range("E11").Value = '  +1.27%  '

$ws.Range("D12").Value = '1.871.01'
$ws.Range("E12").Value = '  +4.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07339'
$ws.Range("E13").Value = '  +3.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.452'
$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.06'
$ws.Range("E15").Value = '  +2.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.588'
$ws.Range("E16").Value = '  +1.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.010'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008814'
$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("E19").Value = '  +0.85%  '

$ws.Range("D20").Value = '27.703.00'
$ws.Range("E20").Value = '  +2.92%  '

$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("D24").Value = '2.116.55'
$ws.Range("E24").Value = '  +4.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.908'
$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.89'
$ws.Range("E26").Value = '  +0.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.63'
$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.144'
$ws.Range("E28").Value = '  +1.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.242'
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.41'
$ws.Range("E30").Value = '  +1.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09003'
$ws.Range("E31").Value = '  +1.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7571'
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.180'
$ws.Range("E33").Value = '  +1.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.548'
$ws.Range("E34").Value = '  +1.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.946'
$ws.Range("E35").Value = '  +1.38%  '

$ws.Range("E36").Value = '  +0.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  +1.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05357'
$ws.Range("E38").Value = '  +1.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01956'
$ws.Range("E39").Value = '  +0.51%  '

$ws.Range("E40").Value = '  +0.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.324'
$ws.Range("E41").Value = '  +1.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.399'
$ws.Range("E42").Value = '  +3.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5322'
$ws.Range("E43").Value = '  +0.43%  '

$ws.Range("E44").Value = '  +0.64%  '

$ws.Range("E45").Value = '  +1.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4914'
$ws.Range("E46").Value = '  +1.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.58'
$ws.Range("E47").Value = '  +2.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.15'
$ws.Range("E48").Value = '  +2.19%  '

$ws.Range("E50").Value = '  +1.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06303'
$ws.Range("E51").Value = '  +0.21%  '
